$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# The table "Condicion_Pacientes" currently covers A1:F76. Add a new row
# to the table so it (and its autofilter) expands to A1:F77.
$lo = $ws.ListObjects.Item("Condicion_Pacientes")
$lo.ListRows.Add() | Out-Null

# Copy the formatting (number format / style) of the last existing data
# row down into the newly created row 77.
$ws.Range("A76:F76").Copy()
$ws.Range("A77:F77").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Fill in the new day's data (2020-05-28 => serial date 43979).
$ws.Cells.Item(77, 1).Value = 43979
$ws.Cells.Item(77, 2).Value = 544
$ws.Cells.Item(77, 3).Value = 112
$ws.Cells.Item(77, 4).Value = 473
$ws.Cells.Item(77, 5).Value = 17
$ws.Cells.Item(77, 6).Value = 39

# Match the saved selection state from the edited workbook.
$ws.Range("B77").Select()
